# New crime data collected — weekly CompStat (112th Precinct) refresh.
# Updates the "Volume/Number" + "Report Covering the Week" header strings
# and the Week-to-Date / 28-Day / YTD / 2-Year crime-complaint figures
# (rows 15-28) to reflect the newly collected week's numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header strings (rich-text cells holding the issue number and the date
# range covered by this report).
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  6"
$ws.Range("C9").Value = "Report Covering the Week  2/3/2025  Through  2/9/2025"

# ---------------------------------------------------------------------
# Helpers for the handful of cells that flip between a numeric value and
# one of the two textual placeholders ("0" for a not-applicable count,
# "***.*" for a not-applicable percent-change) used throughout the grid.
# Donor cells C14 (placeholder "0") and E14 (placeholder "***.*") keep
# their own value/format untouched by row 14, which this edit does not
# otherwise touch.
# ---------------------------------------------------------------------
function Set-PlaceholderZero($addr) {
    $ws.Range("C14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("C14").Copy()
    $ws.Range($addr).PasteSpecial(-4163)   # xlPasteValues
}

function Set-PlaceholderStar($addr) {
    $ws.Range("E14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("E14").Copy()
    $ws.Range($addr).PasteSpecial(-4163)   # xlPasteValues
}

function Set-Num1($addr, $val) {
    # one-decimal numeric style (matches existing style used by K/L/M/N cols)
    $ws.Range($addr).NumberFormat = "#,##0.0;""-""#,##0.0"
    $ws.Range($addr).Value = $val
}

function Set-Int($addr, $val) {
    # plain integer-count style (matches existing style used by C..J cols)
    $ws.Range($addr).NumberFormat = "#,##0"
    $ws.Range($addr).Value = $val
}

# ---------------------------------------------------------------------
# Row 15 — Murder
# ---------------------------------------------------------------------
Set-Num1 "L15" 100

# ---------------------------------------------------------------------
# Row 16 — Rape
# ---------------------------------------------------------------------
Set-PlaceholderZero "C16"
Set-Int "D16" 1
Set-Num1 "E16" -100
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 13
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = 85.714285714285
$ws.Range("L16").Value = 85.714285714285
$ws.Range("M16").Value = -43.478260869565
$ws.Range("N16").Value = -86.868686868686

# ---------------------------------------------------------------------
# Row 17 — Robbery
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 3
Set-Int "D17" 5
Set-Num1 "E17" -40
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 200
$ws.Range("I17").Value = 18
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 80
$ws.Range("L17").Value = 125
$ws.Range("M17").Value = 125
$ws.Range("N17").Value = 125

# ---------------------------------------------------------------------
# Row 18 — Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 340
$ws.Range("I18").Value = 27
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = 145.454545454545
$ws.Range("L18").Value = 12.5
$ws.Range("M18").Value = 107.692307692308
$ws.Range("N18").Value = -83.832335329341

# ---------------------------------------------------------------------
# Row 19 — Burglary
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -10
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 15.384615384615
$ws.Range("I19").Value = 43
$ws.Range("J19").Value = 39
$ws.Range("K19").Value = 10.256410256410
$ws.Range("L19").Value = -6.521739130434
$ws.Range("M19").Value = -18.867924528301
$ws.Range("N19").Value = -53.763440860215

# ---------------------------------------------------------------------
# Row 20 — Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 16
$ws.Range("J20").Value = 10
$ws.Range("K20").Value = 60
$ws.Range("L20").Value = -15.789473684210
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = -96.529284164859

# ---------------------------------------------------------------------
# Row 21 — G.L.A. (bold TOTAL-style row)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -5.555555555555
$ws.Range("F21").Value = 86
$ws.Range("G21").Value = 51
$ws.Range("H21").Value = 68.627450980392
$ws.Range("I21").Value = 119
$ws.Range("J21").Value = 78
$ws.Range("K21").Value = 52.564102564102
$ws.Range("L21").Value = 13.333333333333
$ws.Range("M21").Value = 5.309734513274
$ws.Range("N21").Value = -85.662650602409

# ---------------------------------------------------------------------
# Row 22 — TOTAL
# ---------------------------------------------------------------------
Set-Int "C22" 2
Set-PlaceholderZero "D22"
Set-PlaceholderStar "E22"
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 33.333333333333

# ---------------------------------------------------------------------
# Row 23 — Transit
# ---------------------------------------------------------------------
Set-Num1 "L23" -100

# ---------------------------------------------------------------------
# Row 24 — Housing
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 50
$ws.Range("D24").Value = 34
$ws.Range("E24").Value = 47.058823529411
$ws.Range("F24").Value = 161
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = 21.969696969697
$ws.Range("I24").Value = 220
$ws.Range("J24").Value = 170
$ws.Range("K24").Value = 29.411764705882
$ws.Range("L24").Value = 20.879120879120
$ws.Range("M24").Value = 150

# ---------------------------------------------------------------------
# Row 25 — Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 36
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 56.521739130434
$ws.Range("F25").Value = 132
$ws.Range("G25").Value = 97
$ws.Range("H25").Value = 36.082474226804
$ws.Range("I25").Value = 180
$ws.Range("J25").Value = 128
$ws.Range("K25").Value = 40.625
$ws.Range("L25").Value = 34.328358208955

# ---------------------------------------------------------------------
# Row 26 — Retail Theft
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 133.333333333333
$ws.Range("F26").Value = 31
$ws.Range("H26").Value = 121.428571428571
$ws.Range("I26").Value = 37
$ws.Range("J26").Value = 19
$ws.Range("K26").Value = 94.736842105263
$ws.Range("L26").Value = 68.181818181818
$ws.Range("M26").Value = 42.307692307692

# ---------------------------------------------------------------------
# Row 27 — Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("L27").Value = 50

# ---------------------------------------------------------------------
# Row 28 — UCR Rape*
# ---------------------------------------------------------------------
Set-Int "C28" 1
Set-PlaceholderZero "D28"
Set-PlaceholderStar "E28"
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 300
$ws.Range("I28").Value = 4
$ws.Range("K28").Value = 33.333333333333
